# Updates cryptos price list cells per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string into a given cell as literal TEXT
# (matches source data where prices like "215.58" are stored as text, not numbers),
# then strip the Text number-format back off the cell via a format-only paste from
# a neighboring default-formatted cell so the cell keeps its original (default) style.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range("B2").Copy()
    $cell.PasteSpecial(-4122)
}

$ws.Range("D2").Value = "27.205.56"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "1.682.47"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "215.58"
$ws.Range("E5").Value = "  +0.28%  "
Set-TextValue $ws.Range("D6") "0.518"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +2.01%  "
Set-TextValue $ws.Range("D9") "21.61"
$ws.Range("E9").Value = "  +6.29%  "
Set-TextValue $ws.Range("D10") "0.0624"
$ws.Range("E10").Value = "  +0.57%  "
Set-TextValue $ws.Range("D11") "0.0889"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.919.32"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.680.95"
$ws.Range("E13").Value = "  +0.35%  "
Set-TextValue $ws.Range("D14") "4.16"
$ws.Range("E14").Value = "  +1.71%  "
Set-TextValue $ws.Range("D15") "0.542"
$ws.Range("E15").Value = "  +2.55%  "
Set-TextValue $ws.Range("D16") "66.41"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "27.179.28"
$ws.Range("E17").Value = "  +0.71%  "
Set-TextValue $ws.Range("D18") "238.61"
$ws.Range("E18").Value = "  +0.59%  "
Set-TextValue $ws.Range("D19") "8.08"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "0.0₃0744"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  +0.06%  "
Set-TextValue $ws.Range("D22") "4.55"
$ws.Range("E22").Value = "  +2.58%  "
Set-TextValue $ws.Range("D23") "9.49"
$ws.Range("E23").Value = "  +3.11%  "
Set-TextValue $ws.Range("D24") "2.11"
$ws.Range("E24").Value = "  -3.54%  "
Set-TextValue $ws.Range("D25") "148.18"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  +0.16%  "
Set-TextValue $ws.Range("D27") "16.34"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  +0.61%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue $ws.Range("D30") "0.0499"
$ws.Range("E30").Value = "  +0.34%  "
Set-TextValue $ws.Range("D31") "1.17"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "1.572.22"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("E33").Value = "  +1.48%  "
Set-TextValue $ws.Range("D34") "3.23"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "0.939"
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D38") "2.39"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("E40").Value = "  +3.74%  "
Set-TextValue $ws.Range("D41") "69.15"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "1.827.80"
$ws.Range("E45").Value = "  +0.55%  "
Set-TextValue $ws.Range("D46") "0.788"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.104"
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "8.15"
$ws.Range("E51").Value = "  +5.95%  "
